$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits for existing rows 20-26 (before the new row is inserted) ---

# Rows 21-23: EARNED (VL) credited 1.25 each (G mirrors EARNED via formula)
$ws.Range("C21").Value = 1.25
$ws.Range("C22").Value = 1.25
$ws.Range("C23").Value = 1.25

# Row 24: a Vacation Leave (1 day) taken, noted "VL(1-0-0)", approved date in K
$ws.Range("B24").Value = "VL(1-0-0)"
$ws.Range("C24").Value = 1.25
$ws.Range("D24").Value = 1
$ws.Range("K10").Copy()
$ws.Range("K24").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K24").Value = 45233

# Row 25: EARNED credited 1.25
$ws.Range("C25").Value = 1.25

# --- Insert a new "year separator" row 27 for 2024, pushing the Jan-2024+ rows down ---
$ws.Rows("27").Insert()

# Restore formatting of the newly inserted row 27 from the row above it
$ws.Range("A26:K26").Copy()
$ws.Range("A27:K27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the EARNED formula wiped out by PasteSpecial(formats)
$earnedFormula = 'IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G27").Formula = "=" + $earnedFormula

# Give the new row 27 the "year label" look (same style already used by the 2022/2023 rows)
$ws.Range("A14").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A27").Value = "'2024"

# --- Table1 needs to grow by the inserted row to keep covering the whole data range ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K133"))
$ws.Range("G133").Formula = "=" + $earnedFormula

# --- Remaining data edits, now on the ORIGINAL row 26 (still row 26; insert only shifted rows >= 27) ---
$ws.Range("B26").Value = "VL(3-0-0)"
$ws.Range("D26").Value = 3
$ws.Range("K26").Value = "12/27-29/2023"
